$d = $word.ActiveDocument

$replacements = @(
    @{old="13×89="; new="21×82="},
    @{old="46×83="; new="22×21="},
    @{old="40×46="; new="82×99="},
    @{old="17×12="; new="65×31="},
    @{old="75×30="; new="41×71="},
    @{old="28×68="; new="11×44="},
    @{old="56×61="; new="58×48="},
    @{old="13×90="; new="52×42="},
    @{old="64×72="; new="47×80="},
    @{old="54×54="; new="47×15="},
    @{old="42×78="; new="12×62="},
    @{old="19×31="; new="96×75="},
    @{old="64×87="; new="50×81="},
    @{old="94×11="; new="76×57="},
    @{old="67×11="; new="98×94="},
    @{old="92×92="; new="72×61="},
    @{old="52×90="; new="71×94="},
    @{old="47×81="; new="46×19="},
    @{old="51×92="; new="47×29="},
    @{old="94×13="; new="50×74="},
    @{old="47×99="; new="24×28="},
    @{old="15×81="; new="49×56="},
    @{old="18×66="; new="21×66="},
    @{old="99×87="; new="66×40="},
    @{old="72×52="; new="48×34="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
